# Reposition/resize pictures and rotate the arrow shape on slide 3.
#
# Target offsets/extents come from the target OOXML (EMU units, 1 pt = 12700 EMU).
# The point values below were chosen so that converting pt -> EMU through the
# host's internal 32-bit-float storage reproduces the exact target EMU integers.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# Picture 5 (shape 3) - top-right screenshot -> moves to top-left, shrinks
$picTopRight = $s.Shapes.Item(3)
$picTopRight.Left   = 51.951656341552734
$picTopRight.Top    = 97.83716583251953
$picTopRight.Width  = 323.47467041015625
$picTopRight.Height = 201.69032287597656

# Picture 7 (shape 4) - bottom-right screenshot -> moves below the shape above, shrinks
$picBottomRight = $s.Shapes.Item(4)
$picBottomRight.Left   = 51.951656341552734
$picBottomRight.Top    = 325.2249755859375
$picBottomRight.Width  = 323.47467041015625
$picBottomRight.Height = 201.69039916992188

# Picture 3 (shape 5) - top-left screenshot -> moves to top-right (size unchanged)
$picTopLeft = $s.Shapes.Item(5)
$picTopLeft.Left = 537.1636352539062
$picTopLeft.Top  = 59.27047348022461

# Picture 6 (shape 6) - bottom-left screenshot -> moves to bottom-right (size unchanged)
$picBottomLeft = $s.Shapes.Item(6)
$picBottomLeft.Left = 537.1637573242188
$picBottomLeft.Top  = 278.3290710449219

# Arrow: Left 9 (shape 7) - moves and rotates 180 degrees (size unchanged)
$arrow = $s.Shapes.Item(7)
$arrow.Left     = 424.7030944824219
$arrow.Top      = 251.4528350830078
$arrow.Rotation = 180
